# Refresh the cryptos list with the latest scraped prices/volumes.
# Note: several "Price" values look like plain numbers (e.g. 1.00, 0.999,
# 17.71 ...). Excel's Range.Value setter auto-coerces such strings to
# numeric cells, which would silently drop meaningful trailing zeros
# (e.g. "1.00" -> 1). To keep these as text - matching how the source
# data is stored - they are written with a leading apostrophe, exactly
# like a user typing a text-forced value directly into a cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '98.005.24'
$ws.Range('E2').Value = '  -0.94%  '
$ws.Range('D3').Value = '3.434.47'
$ws.Range('E3').Value = '  +4.03%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = '''256.60'
$ws.Range('E5').Value = '  +0.84%  '
$ws.Range('D6').Value = '''657.61'
$ws.Range('E6').Value = '  +5.44%  '
$ws.Range('D7').Value = '''1.48'
$ws.Range('E7').Value = '  +1.75%  '
$ws.Range('D8').Value = '''0.432'
$ws.Range('E8').Value = '  +6.06%  '
$ws.Range('D9').Value = '''1.06'
$ws.Range('E9').Value = '  +10.20%  '
$ws.Range('E10').Value = '  -0.13%  '
$ws.Range('D11').Value = '3.431.45'
$ws.Range('E11').Value = '  +4.04%  '
$ws.Range('E12').Value = '  +6.81%  '
$ws.Range('D13').Value = '''42.37'
$ws.Range('E13').Value = '  +6.93%  '
$ws.Range('D14').Value = '''6.57'
$ws.Range('E14').Value = '  +20.05%  '
$ws.Range('D15').Value = '''0.0000260'
$ws.Range('E15').Value = '  +4.20%  '
$ws.Range('D16').Value = '97.784.95'
$ws.Range('E16').Value = '  -1.06%  '
$ws.Range('D17').Value = '4.068.70'
$ws.Range('E17').Value = '  +3.77%  '
$ws.Range('D18').Value = '''8.71'
$ws.Range('E18').Value = '  +38.71%  '
$ws.Range('D19').Value = '3.425.41'
$ws.Range('E19').Value = '  +3.62%  '
$ws.Range('D20').Value = '''17.71'
$ws.Range('E20').Value = '  +14.61%  '
$ws.Range('D21').Value = '''0.525'
$ws.Range('E21').Value = '  +71.32%  '
$ws.Range('D22').Value = '''10.97'
$ws.Range('E22').Value = '  +16.68%  '
$ws.Range('E23').Value = '  +0.51%  '
$ws.Range('D24').Value = '''511.12'
$ws.Range('E24').Value = '  +5.51%  '
$ws.Range('D25').Value = '''0.0000207'
$ws.Range('E25').Value = '  +2.48%  '
$ws.Range('B26').Value = 'NEARProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D26').Value = '''6.18'
$ws.Range('E26').Value = '  +10.07%  '
$ws.Range('B27').Value = 'Litecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D27').Value = '''99.18'
$ws.Range('E27').Value = '  +11.62%  '
$ws.Range('D28').Value = '''12.78'
$ws.Range('E28').Value = '  +7.15%  '
$ws.Range('D29').Value = '''0.153'
$ws.Range('E29').Value = '  +12.89%  '
$ws.Range('D30').Value = '''11.50'
$ws.Range('E30').Value = '  +12.65%  '
$ws.Range('D31').Value = '''1.00'
$ws.Range('E31').Value = '  +0.10%  '
$ws.Range('D32').Value = '''0.197'
$ws.Range('E32').Value = '  +4.85%  '
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('D34').Value = '''0.574'
$ws.Range('E34').Value = '  +22.17%  '
$ws.Range('D35').Value = '''29.99'
$ws.Range('E35').Value = '  +8.06%  '
$ws.Range('E36').Value = '  +13.49%  '
$ws.Range('D37').Value = '''7.90'
$ws.Range('E37').Value = '  +9.81%  '
$ws.Range('D38').Value = '''0.158'
$ws.Range('E38').Value = '  +7.01%  '
$ws.Range('D39').Value = '''1.41'
$ws.Range('E39').Value = '  +14.76%  '
$ws.Range('D40').Value = '''516.41'
$ws.Range('E40').Value = '  +5.83%  '
$ws.Range('D41').Value = '''24.72'
$ws.Range('E41').Value = '  -0.36%  '
$ws.Range('D42').Value = '''0.857'
$ws.Range('E42').Value = '  +9.74%  '
$ws.Range('D43').Value = '''0.0421'
$ws.Range('E43').Value = '  +27.18%  '
$ws.Range('D44').Value = '''3.68'
$ws.Range('E44').Value = '  +1.38%  '
$ws.Range('D45').Value = '''3.33'
$ws.Range('E45').Value = '  +7.11%  '
$ws.Range('D46').Value = '''5.45'
$ws.Range('E46').Value = '  +16.06%  '
$ws.Range('B47').Value = 'Cosmos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D47').Value = '''8.23'
$ws.Range('E47').Value = '  +12.99%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').Value = '''1.00'
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('D49').Value = '''1.60'
$ws.Range('E49').Value = '  +18.04%  '
$ws.Range('D50').Value = '''2.11'
$ws.Range('E50').Value = '  +8.46%  '
$ws.Range('D51').Value = '''50.97'
$ws.Range('E51').Value = '  +10.99%  '
